# Fix bug: quarterly history imported as monthly for both DoEstimationwithData and xxSV.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 0.05
$ws.Range("F2").Value = 0.03

# Row 3
$ws.Range("E3").Value = 0.05
$ws.Range("F3").Value = 0.03
$ws.Range("J3").Value = 1

# Row 4
$ws.Range("E4").Value = 0.05
$ws.Range("F4").Value = 0.03
$ws.Range("J4").Value = 1

# Row 5
$ws.Range("E5").Value = 0.7
$ws.Range("F5").Value = 0.38
$ws.Range("G5").Value = 0.88
$ws.Range("I5").Value = 0.02
$ws.Range("J5").Value = 0.03
$ws.Range("K5").Value = 0.98

# Row 6
$ws.Range("E6").Value = 0.7
$ws.Range("F6").Value = 0.38
$ws.Range("G6").Value = 0.88
$ws.Range("H6").Value = 0.2
$ws.Range("I6").Value = 0.01
$ws.Range("J6").Value = 0.03
$ws.Range("K6").Value = 0.98
$ws.Range("L6").Value = 1.11
